$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: issue number and week-covering dates ---
$ws.Range("A8").Value = "Volume 33   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/19/2026  Through  1/25/2026"

# --- Weekly crime statistics table (rows 14-28), TOTAL row 21, and related rows ---
# Row 14
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("G14").Value = 1
$ws.Range("J14").NumberFormat = "#,##0"
$ws.Range("J14").Value = 1
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E14").Value = -100
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H14").Value = -100
$ws.Range("K14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K14").Value = -100

# Row 15
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 2
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 4
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 1
$ws.Range("J15").NumberFormat = "#,##0"
$ws.Range("J15").Value = 4
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -50
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H15").Value = -75
$ws.Range("K15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K15").Value = -75
$ws.Range("L15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L15").Value = -50
$ws.Range("M15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M15").Value = -50
$ws.Range("N15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N15").Value = -75

# Row 16
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 6
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 6
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("F16").Value = 22
$ws.Range("G16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 19
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("I16").Value = 21
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("J16").Value = 13
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E16").Value = 0
$ws.Range("H16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H16").Value = 15.78947368421
$ws.Range("K16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K16").Value = 61.538461538461
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L16").Value = -16
$ws.Range("M16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M16").Value = 16.666666666666
$ws.Range("N16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N16").Value = -78.571428571428

# Row 17
$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("F17").Value = 24
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 43
$ws.Range("I17").NumberFormat = "#,##0"
$ws.Range("I17").Value = 19
$ws.Range("J17").NumberFormat = "#,##0"
$ws.Range("J17").Value = 38
$ws.Range("H17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H17").Value = -44.186046511627
$ws.Range("K17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K17").Value = -50
$ws.Range("L17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L17").Value = -20.833333333333
$ws.Range("M17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M17").Value = 18.75
$ws.Range("N17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N17").Value = -40.625

# Row 18
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 1
$ws.Range("F18").NumberFormat = "#,##0"
$ws.Range("F18").Value = 2
$ws.Range("G18").NumberFormat = "#,##0"
$ws.Range("G18").Value = 7
$ws.Range("J18").NumberFormat = "#,##0"
$ws.Range("J18").Value = 6
$ws.Range("H18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H18").Value = -71.428571428571
$ws.Range("K18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K18").Value = -83.333333333333
$ws.Range("L18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L18").Value = -95.454545454545
$ws.Range("M18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M18").Value = -95
$ws.Range("N18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N18").Value = -99.532710280373

# Row 19
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("C19").Value = 9
$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("D19").Value = 10
$ws.Range("F19").NumberFormat = "#,##0"
$ws.Range("F19").Value = 44
$ws.Range("G19").NumberFormat = "#,##0"
$ws.Range("G19").Value = 48
$ws.Range("I19").NumberFormat = "#,##0"
$ws.Range("I19").Value = 40
$ws.Range("J19").NumberFormat = "#,##0"
$ws.Range("J19").Value = 43
$ws.Range("E19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E19").Value = -10
$ws.Range("H19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H19").Value = -8.333333333333
$ws.Range("K19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K19").Value = -6.976744186046
$ws.Range("L19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L19").Value = -50
$ws.Range("M19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M19").Value = -2.439024390243
$ws.Range("N19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N19").Value = -54.545454545454

# Row 20
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 5
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 4
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("F20").Value = 18
$ws.Range("G20").NumberFormat = "#,##0"
$ws.Range("G20").Value = 12
$ws.Range("I20").NumberFormat = "#,##0"
$ws.Range("I20").Value = 14
$ws.Range("J20").NumberFormat = "#,##0"
$ws.Range("J20").Value = 10
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E20").Value = 25
$ws.Range("H20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H20").Value = 50
$ws.Range("K20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K20").Value = 40
$ws.Range("L20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L20").Value = -30
$ws.Range("M20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M20").Value = -17.647058823529
$ws.Range("N20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N20").Value = -90.666666666666

# Row 21
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("C21").Value = 25
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 31
$ws.Range("F21").NumberFormat = "#,##0"
$ws.Range("F21").Value = 111
$ws.Range("G21").NumberFormat = "#,##0"
$ws.Range("G21").Value = 134
$ws.Range("I21").NumberFormat = "#,##0"
$ws.Range("I21").Value = 96
$ws.Range("J21").NumberFormat = "#,##0"
$ws.Range("J21").Value = 115
$ws.Range("E21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("E21").Value = -19.354838709677
$ws.Range("H21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("H21").Value = -17.164179104477
$ws.Range("K21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("K21").Value = -16.521739130434
$ws.Range("L21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("L21").Value = -44.508670520231
$ws.Range("M21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("M21").Value = -15.78947368421
$ws.Range("N21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("N21").Value = -83.645655877342

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 5
$ws.Range("J22").NumberFormat = "#,##0"
$ws.Range("J22").Value = 5
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H22").Value = -40
$ws.Range("K22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K22").Value = -60
$ws.Range("L22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L22").Value = -60
$ws.Range("M22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M22").Value = 0

# Row 24
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("C24").Value = 20
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("D24").Value = 32
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("F24").Value = 92
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("G24").Value = 99
$ws.Range("I24").NumberFormat = "#,##0"
$ws.Range("I24").Value = 79
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("J24").Value = 86
$ws.Range("E24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E24").Value = -37.5
$ws.Range("H24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H24").Value = -7.070707070707
$ws.Range("K24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K24").Value = -8.13953488372
$ws.Range("L24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L24").Value = -51.829268292682
$ws.Range("M24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M24").Value = 16.176470588235

# Row 25
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("C25").Value = 5
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 17
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("F25").Value = 19
$ws.Range("G25").NumberFormat = "#,##0"
$ws.Range("G25").Value = 39
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("I25").Value = 17
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("J25").Value = 38
$ws.Range("E25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E25").Value = -70.588235294117
$ws.Range("H25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H25").Value = -51.282051282051
$ws.Range("K25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K25").Value = -55.263157894736
$ws.Range("L25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L25").Value = -83

# Row 26
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 11
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 14
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 68
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 76
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("I26").Value = 62
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 69
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = -21.428571428571
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H26").Value = -10.526315789473
$ws.Range("K26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K26").Value = -10.144927536231
$ws.Range("L26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L26").Value = 1.639344262295
$ws.Range("M26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M26").Value = 1.639344262295

# Row 27
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 2
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 2
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 5
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 2
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("J27").Value = 5
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -50
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H27").Value = -60
$ws.Range("K27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K27").Value = -60
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L27").Value = -71.428571428571

# Row 28
$ws.Range("F28").NumberFormat = "#,##0"
$ws.Range("F28").Value = 4
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("G28").Value = 6
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("J28").Value = 5
$ws.Range("H28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H28").Value = -33.333333333333
$ws.Range("K28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K28").Value = -40
$ws.Range("L28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L28").Value = -50

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"

